# Remove the trailing ";" from the Python code-sample lines shown on the
# "input/output de Dados" slides (the statements aren't terminated by a
# semicolon in Python). Each target run is located with
# TextFrame.TextRange.Characters(start, length) - computed from the shape's
# original OOXML - and rewritten in place so only that run's text changes
# while every other run keeps its original formatting (rPr) untouched.
# Because every edit removes exactly one character, edits within the same
# shape are applied from the highest Start offset down to the lowest so
# that earlier (still-pending) offsets are not invalidated.
$p = $ppt.ActivePresentation

$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(37, 8).Text = '(x + y)'
$tr.Characters(22, 8).Text = [char]0x9 + 'y = 43'
$tr.Characters(14, 7).Text = 'x = 15'
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(41, 11).Text = '(x + y) # '
$tr.Characters(21, 13).Text = [char]0x9 + 'y = "Julia' + [char]0x201C
$tr.Characters(14, 6).Text = 'x = 8'
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(41, 10).Text = '(x, y) # '
$tr.Characters(21, 13).Text = [char]0x9 + 'y = "Julia' + [char]0x201C
$tr.Characters(14, 6).Text = 'x = 8'
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(141, 9).Text = '(valor))'
$tr.Characters(92, 12).Text = [char]0x9 + 'valor = 35'
$tr.Characters(70, 9).Text = '(valor))'
$tr.Characters(39, 13).Text = ' ' + [char]0xE9 + ' {} reais' + [char]0x201C
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(65, 9).Text = '(valor))'
$tr.Characters(12, 14).Text = [char]0x9 + 'valor = 35.6'
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(159, 3).Text = '))'
$tr.Characters(103, 16).Text = ' {:.2f} reais.' + [char]0x201C
$tr.Characters(46, 6).Text = ' = 34'
$tr.Characters(32, 7).Text = ' = 135'
$tr.Characters(14, 10).Text = [char]0x9 + 'qtd = 12'
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(204, 3).Text = '))'
$tr.Characters(105, 20).Text = ' {2:.2f} reais.' + [char]0x2019 + ' # '
$tr.Characters(46, 6).Text = ' = 34'
$tr.Characters(32, 7).Text = ' = 135'
$tr.Characters(14, 10).Text = [char]0x9 + 'qtd = 12'
$s = $p.Slides.Item(17)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(115, 3).Text = '))'
$tr.Characters(82, 3).Text = '.' + [char]0x201C
$tr.Characters(30, 11).Text = ' = "Julia' + [char]0x201C
$tr.Characters(18, 6).Text = ' = 25'
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(107, 11).Text = ' = ' + [char]0x2018 + 'Gol' + [char]0x2019 + '))'
$tr.Characters(58, 4).Text = '}.' + [char]0x2019
$s = $p.Slides.Item(19)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(63, 4).Text = '}")'
$tr.Characters(39, 2).Text = [char]0x201C
$tr.Characters(19, 6).Text = ' = 17'
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(111, 21).Text = '("Hello, World!") # '
$tr.Characters(86, 18).Text = '("Hello, World!")'
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(98, 24).Text = '("Meu nome ' + [char]0xE9 + ': " + nome)'
$tr.Characters(72, 20).Text = '("Informe o nome:")'
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(125, 26).Text = '("Meu nome ' + [char]0xE9 + ': ' + [char]0x201C + ', salario)'
$tr.Characters(95, 24).Text = '("Informe o sal' + [char]0xE1 + 'rio:"))'
$tr.Characters(48, 24).Text = '("Meu nome ' + [char]0xE9 + ': ' + [char]0x201C + ', idade)'
$tr.Characters(18, 24).Text = '("Informe sua idade:"))'
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(45, 4).Text = '(x)'
$tr.Characters(14, 24).Text = 'x = "Python ' + [char]0xE9 + ' incr' + [char]0xED + 'vel' + [char]0x201C
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(61, 10).Text = '(x, y, z)'
$tr.Characters(38, 16).Text = [char]0x9 + 'z = ' + [char]0x201C + 'incr' + [char]0xED + 'vel' + [char]0x201C
$tr.Characters(28, 9).Text = [char]0x9 + 'y = ' + [char]0x201C + [char]0xE9 + [char]0x201C
$tr.Characters(14, 13).Text = 'x = "Python' + [char]0x201C
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(61, 12).Text = '(x + y + z)'
$tr.Characters(38, 16).Text = [char]0x9 + 'z = ' + [char]0x201C + 'incr' + [char]0xED + 'vel' + [char]0x201C
$tr.Characters(28, 9).Text = [char]0x9 + 'y = ' + [char]0x201C + [char]0xE9 + [char]0x201C
$tr.Characters(14, 13).Text = 'x = "Python' + [char]0x201C

